# Updated cryptos list on Wed Jan 10 02:49:28 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns for every coin row,
# and swaps row 51 from "EnergySwap" to "ordi" (name, link, price, volume).
#
# The Price/Link/Coin columns are plain-text cells (General number format)
# holding strings that often *look* numeric ("301.23", "1.00", "0.999", ...).
# Writing such a string straight into Range.Value would make Excel silently
# reinterpret it as a real number (dropping significant trailing zeros such
# as "2.20" -> 2.2). Set-CellText works around that the same way a person
# would in the Excel UI: force the cell to Text (@) before writing, then
# restore the original "Normal" style so no other formatting is disturbed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $looksNumeric = $text -match '^[+-]?[0-9]+(\.[0-9]+)?$'

    if ($looksNumeric) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

# row, Price(D), Volume1h(E)
$rows = @(
    @(2,  "45.984.75",  "  -1.62%  "),
    @(3,  "2.352.25",   "  +1.63%  "),
    @(4,  $null,        "  +0.00%  "),
    @(5,  "301.23",     "  +0.34%  "),
    @(6,  "98.84",      "  +0.04%  "),
    @(7,  "0.568",      "  -0.97%  "),
    @(8,  $null,        "  +0.04%  "),
    @(9,  "0.512",      "  -3.57%  "),
    @(10, "34.44",      "  -3.20%  "),
    @(11, "0.0797",     "  -0.22%  "),
    @(12, "7.12",       "  -3.29%  "),
    @(13, $null,        "  -0.35%  "),
    @(14, "2.710.15",   "  +1.66%  "),
    @(15, "2.356.08",   "  +1.75%  "),
    @(16, "0.808",      "  -1.39%  "),
    @(17, "13.59",      "  -2.77%  "),
    @(18, "45.881.24",  "  -1.74%  "),
    @(19, "12.72",      "  -4.78%  "),
    @(20, "0.0₃0967",   "  +2.73%  "),
    @(21, "6.01",       "  -2.33%  "),
    @(22, "67.26",      "  +0.57%  "),
    @(23, "244.48",     "  -1.93%  "),
    @(24, $null,        "  -2.84%  "),
    @(25, "0.999",      "  -0.26%  "),
    @(26, "1.92",       "  -3.59%  "),
    @(27, "39.76",      "  -7.92%  "),
    @(28, "2.20",       "  -1.61%  "),
    @(29, "9.81",       "  -0.57%  "),
    @(30, "3.74",       "  +19.32%  "),
    @(31, "20.82",      "  +3.44%  "),
    @(32, $null,        "  +5.20%  "),
    @(33, "5.52",       "  -4.84%  "),
    @(34, "146.08",     "  -1.10%  "),
    @(35, "0.0774",     "  -3.37%  "),
    @(36, "0.113",      "  +0.62%  "),
    @(37, "1.88",       "  +3.70%  "),
    @(38, $null,        "  -2.27%  "),
    @(39, "15.09",      "  -3.46%  "),
    @(40, "3.94",       "  -1.54%  "),
    @(41, "0.0301",     "  -2.35%  "),
    @(42, "3.21",       "  -7.48%  "),
    @(43, "1.878.65",   "  +1.86%  "),
    @(44, "1.00",       "  +0.04%  "),
    @(45, "91.46",      "  +0.96%  "),
    @(46, "1.78",       "  -10.89%  "),
    @(47, "0.186",      "  -7.04%  "),
    @(48, "8.27",       "  +3.69%  "),
    @(49, "97.69",      "  +0.21%  "),
    @(50, "2.581.08",   "  +1.44%  ")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $price  = $r[1]
    $volume = $r[2]
    if ($null -ne $price) {
        Set-CellText $rowNum "D" $price
    }
    Set-CellText $rowNum "E" $volume
}

# Row 51: EnergySwap -> ordi (name, link, price, volume all change)
$ws.Cells.Item(51, "B").Value = "ordi"
$ws.Cells.Item(51, "C").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
Set-CellText 51 "D" "68.47"
Set-CellText 51 "E" "  -9.52%  "
